$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 64 (this shifts the existing rows 64..106 down to 65..107,
# carrying the old row 64 data down to the new row 65, and the previous last
# data row (old 106) down to the new last row 107).
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with a fresh weekly price record.
# Columns that repeat the same record metadata as the row above (market,
# region, product codes, quality) are carried over explicitly since the
# row-insert leaves them blank.
$ws.Range("A64").Value = 11
$ws.Range("B64").Value = "Vega Monumental Concepción"
$ws.Range("C64").Value = "Bíobío"
$ws.Range("D64").Value = 44603
$ws.Range("E64").Value = 8
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100109
$ws.Range("H64").Value = "Uva"
$ws.Range("I64").Value = 100109001
$ws.Range("J64").Value = "Uva"
$ws.Range("K64").Value = "Red Globe"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 200
$ws.Range("N64").Value = 13000
$ws.Range("O64").Value = 14000
$ws.Range("P64").Value = 13500
$ws.Range("Q64").Value = "$/bandeja 18 kilos"
$ws.Range("R64").Value = "Provincia de Limarí"
$ws.Range("S64").Value = 750
$ws.Range("T64").Value = 18
